$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.020.81'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '1.828.22'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9993'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.52'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6200'
$ws.Range('E6').Value = '  -6.66%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.44'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07440'
$ws.Range('E9').Value = '  +0.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.2913'
$ws.Range('E10').Value = '  -0.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.58'
$ws.Range('E11').Value = '  -0.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07604'
$ws.Range('E12').Value = '  -1.84%  '
$ws.Range('D13').Value = '1.825.63'
$ws.Range('E13').Value = '  -0.54%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.944'
$ws.Range('E14').Value = '  -0.84%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6606'
$ws.Range('E15').Value = '  -1.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '81.84'
$ws.Range('E16').Value = '  -1.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000009069'
$ws.Range('E17').Value = '  +8.70%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.959'
$ws.Range('E18').Value = '  -2.19%  '
$ws.Range('D19').Value = '29.022.62'
$ws.Range('E19').Value = '  -0.39%  '
$ws.Range('D20').Value = '2.078.11'
$ws.Range('E20').Value = '  -0.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '223.64'
$ws.Range('E21').Value = '  -2.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.31'
$ws.Range('E22').Value = '  -1.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.150'
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.371'
$ws.Range('E27').Value = '  -2.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1350'
$ws.Range('E28').Value = '  -3.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '17.76'
$ws.Range('E29').Value = '  -1.38%  '
$ws.Range('E30').Value = '  -0.84%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.204'
$ws.Range('E31').Value = '  +1.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.021'
$ws.Range('E32').Value = '  -0.28%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.036'
$ws.Range('E33').Value = '  -1.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05210'
$ws.Range('E34').Value = '  -1.91%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.828'
$ws.Range('E35').Value = '  -2.10%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.147'
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.7312'
$ws.Range('E37').Value = '  -2.33%  '
$ws.Range('E38').Value = '  +0.22%  '
$ws.Range('D39').Value = '1.273.39'
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('E40').Value = '  +0.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01774'
$ws.Range('E41').Value = '  -1.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.295'
$ws.Range('E42').Value = '  +6.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8926'
$ws.Range('E43').Value = '  -3.89%  '
$ws.Range('E44').Value = '  +0.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.48'
$ws.Range('E45').Value = '  -0.43%  '
$ws.Range('D46').Value = '1.976.61'
$ws.Range('E46').Value = '  -0.48%  '
$ws.Range('E47').Value = '  -0.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '63.18'
$ws.Range('E48').Value = '  +0.25%  '
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.695'
$ws.Range('E50').Value = '  -3.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3949'
$ws.Range('E51').Value = '  -1.77%  '
